$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''284.39'
$ws.Range("E2").Value = '''2.18%'

$ws.Range("D3").Value = '''28.49'
$ws.Range("E3").Value = '''4.12%'

$ws.Range("D4").Value = '''5.092'
$ws.Range("E4").Value = '''5.14%'

$ws.Range("D5").Value = '''0.06636'
$ws.Range("E5").Value = '''4.21%'

$ws.Range("D6").Value = '''7.295'
$ws.Range("E6").Value = '''3.82%'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''3.386'
$ws.Range("E7").Value = '''1.95%'

$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '''1.352'
$ws.Range("E8").Value = '''3.06%'

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9329'
$ws.Range("E9").Value = '''4.56%'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1564'
$ws.Range("E10").Value = '''2.80%'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '''0.06138'
$ws.Range("E11").Value = '''11.57%'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.07581'
$ws.Range("E12").Value = '''2.52%'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.02880'
$ws.Range("E13").Value = '''-2.19%'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.08947'
$ws.Range("E14").Value = '''-0.23%'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001592'
$ws.Range("E15").Value = '''0.63%'

$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04460'
$ws.Range("E16").Value = '''1.37%'

$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.0006404'
$ws.Range("E17").Value = '''0.44%'

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '''0.006136'
$ws.Range("E18").Value = '''-0.26%'

$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = '''3.475'
$ws.Range("E19").Value = '''-0.05%'

$ws.Range("D20").Value = '''2.251'
$ws.Range("E20").Value = '''0.75%'

$ws.Range("D21").Value = '''0.3195'
$ws.Range("E21").Value = '''0.78%'

$ws.Range("D22").Value = '''0.1302'
$ws.Range("E22").Value = '''-3.44%'

$ws.Range("D23").Value = '''4.079'
$ws.Range("E23").Value = '''4.30%'

$ws.Range("D24").Value = '''0.1523'
$ws.Range("E24").Value = '''1.17%'

$ws.Range("D25").Value = '''0.001178'
$ws.Range("E25").Value = '''0.25%'

$ws.Range("D26").Value = '''0.004458'
$ws.Range("E26").Value = '''4.27%'

$ws.Range("D27").Value = '''0.0001246'
$ws.Range("E27").Value = '''5.70%'

$ws.Range("D28").Value = '''0.0001613'
$ws.Range("E28").Value = '''-9.08%'

$ws.Range("D40").Value = '''0.04162'
$ws.Range("E40").Value = '''3.34%'

$ws.Range("D41").Value = '''0.006737'
$ws.Range("E41").Value = '''0.31%'

$ws.Range("D42").Value = '''0.1244'
$ws.Range("E42").Value = '''-11.02%'

$ws.Range("D43").Value = '''0.002013'
$ws.Range("E43").Value = '''-2.63%'

$ws.Range("D44").Value = '''0.01149'
$ws.Range("E44").Value = '''2.87%'

$ws.Range("D45").Value = '''0.00005705'
$ws.Range("E45").Value = '''2.73%'

$ws.Range("E46").Value = '''20.74%'

$ws.Range("D47").Value = '''0.01303'
$ws.Range("E47").Value = '''-29.50%'
